$p = $ppt.ActivePresentation

# --- 1. Change the table style on slide 6 (the "SOURCES OF FINANCE" table) ---
$s = $p.Slides.Item(6)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{1E2EC4F8-0ADD-4385-AE25-73B3CC5B7C5A}")
    }
}
